$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

# Fill in the three new rows (107-109) with id / text_JP / text values.
# Column A = id, Column B = filter (left blank), Column C = text_JP, Column D = text (English)

$ws.Cells.Item(107, 1).Value = "cwl_log_deduplicate"
$ws.Cells.Item(107, 3).Value = "de-duplicated rows: {0}"
$ws.Cells.Item(107, 4).Value = "de-duplicated rows: {0}"

$ws.Cells.Item(108, 1).Value = "cwl_log_unique_count"
$ws.Cells.Item(108, 3).Value = "{0} row count {1} | unique count {2}"
$ws.Cells.Item(108, 4).Value = "{0} row count {1} | unique count {2}"

$ws.Cells.Item(109, 1).Value = "cwl_log_spatial_gen"
$ws.Cells.Item(109, 3).Value = "instantiating new zone {0} / {1}"
$ws.Cells.Item(109, 4).Value = "instantiating new zone {0} / {1}"

# Update selection to reflect the new active cell / selection range.
$ws.Range("D107:D109").Select()

# Update the workbook window size/position to match the new view state.
$win = $excel.ActiveWindow
$win.Left = 3120
$win.Top = 3105
$win.Width = 28185
$win.Height = 16695
